$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append 25 new name rows (rows 6-30), each a first/last name pair,
# continuing straight on from the existing 5 rows of data.
$ws.Range("A6").Value = "Reuben"
$ws.Range("B6").Value = "Homenick"
$ws.Range("A7").Value = "Merri"
$ws.Range("B7").Value = "Reichel"
$ws.Range("A8").Value = "Collette"
$ws.Range("B8").Value = "Lynch"
$ws.Range("A9").Value = "Jonna"
$ws.Range("B9").Value = "Beier"
$ws.Range("A10").Value = "Wilburn"
$ws.Range("B10").Value = "Franecki"
$ws.Range("A11").Value = "Cyrstal"
$ws.Range("B11").Value = "Kovacek"
$ws.Range("A12").Value = "Robt"
$ws.Range("B12").Value = "Hermiston"
$ws.Range("A13").Value = "Houston"
$ws.Range("B13").Value = "McKenzie"
$ws.Range("A14").Value = "Narcisa"
$ws.Range("B14").Value = "Lebsack"
$ws.Range("A15").Value = "Jon"
$ws.Range("B15").Value = "Fritsch"
$ws.Range("A16").Value = "Michal"
$ws.Range("B16").Value = "Greenholt"
$ws.Range("A17").Value = "Kareem"
$ws.Range("B17").Value = "Hauck"
$ws.Range("A18").Value = "Barrett"
$ws.Range("B18").Value = "Wyman"
$ws.Range("A19").Value = "Robby"
$ws.Range("B19").Value = "Graham"
$ws.Range("A20").Value = "Clarine"
$ws.Range("B20").Value = "Luettgen"
$ws.Range("A21").Value = "Emmaline"
$ws.Range("B21").Value = "Hammes"
$ws.Range("A22").Value = "Martin"
$ws.Range("B22").Value = "Hegmann"
$ws.Range("A23").Value = "Ward"
$ws.Range("B23").Value = "Carroll"
$ws.Range("A24").Value = "Sanford"
$ws.Range("B24").Value = "Lakin"
$ws.Range("A25").Value = "Flavia"
$ws.Range("B25").Value = "Upton"
$ws.Range("A26").Value = "Clinton"
$ws.Range("B26").Value = "Marvin"
$ws.Range("A27").Value = "Rickie"
$ws.Range("B27").Value = "Brekke"
$ws.Range("A28").Value = "Amado"
$ws.Range("B28").Value = "Powlowski"
$ws.Range("A29").Value = "Royal"
$ws.Range("B29").Value = "Windler"
$ws.Range("A30").Value = "Peg"
$ws.Range("B30").Value = "Yost"
